$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "87.330.75"
$ws.Range("E2").Value = "  +7.40%  "

$ws.Range("D3").Value = "3.310.60"
$ws.Range("E3").Value = "  +3.46%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "'217.03"
$ws.Range("E5").Value = "  +3.42%  "

$ws.Range("D6").Value = "'647.09"
$ws.Range("E6").Value = "  +1.37%  "

$ws.Range("D7").Value = "'0.349"
$ws.Range("E7").Value = "  +19.38%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("D9").Value = "'0.603"
$ws.Range("E9").Value = "  +1.34%  "

$ws.Range("D10").Value = "3.313.31"
$ws.Range("E10").Value = "  +3.58%  "

$ws.Range("D11").Value = "'0.585"
$ws.Range("E11").Value = "  -2.29%  "

$ws.Range("D12").Value = "'0.0000267"
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.168"
$ws.Range("E13").Value = "  +1.30%  "

$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'35.30"
$ws.Range("E14").Value = "  +8.90%  "

$ws.Range("D15").Value = "3.924.90"
$ws.Range("E15").Value = "  +3.48%  "

$ws.Range("D16").Value = "'5.48"
$ws.Range("E16").Value = "  +1.05%  "

$ws.Range("D17").Value = "87.352.06"
$ws.Range("E17").Value = "  +7.51%  "

$ws.Range("D18").Value = "3.319.15"
$ws.Range("E18").Value = "  +3.86%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'10.45"
$ws.Range("E19").Value = "  +12.18%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'14.70"
$ws.Range("E20").Value = "  +1.31%  "

$ws.Range("B21").Value = "SuiNetwork"
$ws.Range("C21").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D21").Value = "'3.11"
$ws.Range("E21").Value = "  -1.97%  "

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'456.53"
$ws.Range("E22").Value = "  +2.19%  "

$ws.Range("D23").Value = "'5.53"
$ws.Range("E23").Value = "  +4.12%  "

$ws.Range("D24").Value = "'5.48"
$ws.Range("E24").Value = "  +7.14%  "

$ws.Range("D25").Value = "'12.56"
$ws.Range("E25").Value = "  +10.37%  "

$ws.Range("D26").Value = "3.489.45"
$ws.Range("E26").Value = "  +3.62%  "

$ws.Range("D27").Value = "'78.62"
$ws.Range("E27").Value = "  +1.42%  "

$ws.Range("D28").Value = "'0.201"
$ws.Range("E28").Value = "  +55.54%  "

$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").Value = "'0.0000126"
$ws.Range("E30").Value = "  -2.11%  "

$ws.Range("D31").Value = "'606.48"
$ws.Range("E31").Value = "  +6.14%  "

$ws.Range("D32").Value = "'9.36"
$ws.Range("E32").Value = "  +0.98%  "

$ws.Range("D33").Value = "'1.60"
$ws.Range("E33").Value = "  +5.30%  "

$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("D35").Value = "'2.08"
$ws.Range("E35").Value = "  +1.76%  "

$ws.Range("D36").Value = "'7.15"
$ws.Range("E36").Value = "  +19.00%  "

$ws.Range("D37").Value = "'0.147"
$ws.Range("E37").Value = "  -4.05%  "

$ws.Range("D38").Value = "'23.52"
$ws.Range("E38").Value = "  +1.17%  "

$ws.Range("D39").Value = "'2.16"
$ws.Range("E39").Value = "  +2.62%  "

$ws.Range("D40").Value = "'0.419"
$ws.Range("E40").Value = "  +0.53%  "

$ws.Range("D41").Value = "'21.83"
$ws.Range("E41").Value = "  +4.91%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("D43").Value = "'3.00"
$ws.Range("E43").Value = "  -3.46%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'194.35"
$ws.Range("E44").Value = "  +1.35%  "

$ws.Range("D45").Value = "'158.86"
$ws.Range("E45").Value = "  -0.07%  "

$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").Value = "'1.41"
$ws.Range("E47").Value = "  +3.63%  "

$ws.Range("D48").Value = "'46.21"
$ws.Range("E48").Value = "  +6.94%  "

$ws.Range("D49").Value = "'4.46"
$ws.Range("E49").Value = "  +3.19%  "

$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "'0.666"
$ws.Range("E50").Value = "  +2.87%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").Value = "'0.781"
$ws.Range("E51").Value = "  -1.09%  "
